$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 293, shifting existing rows 293:352 down to 294:353.
$ws.Rows("293:293").Insert()

# Populate the newly inserted row 293 with the new data record.
$ws.Range("A293").Value = 10
$ws.Range("B293").Value = "Vega Modelo de Temuco"
$ws.Range("C293").Value = "La Araucanía"
$ws.Range("D293").Value = 45173
$ws.Range("E293").Value = 9
$ws.Range("F293").Value = 100114007
$ws.Range("G293").Value = "Jengibre"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 160
$ws.Range("K293").Value = 22000
$ws.Range("L293").Value = 24000
$ws.Range("M293").Value = 23000
$ws.Range("N293").Value = "$/caja 13 kilos"
$ws.Range("O293").Value = "Perú"
$ws.Range("P293").Value = 1769
$ws.Range("Q293").Value = 13
$ws.Range("R293").Value = "Hortaliza"

# Match the date style used by the rest of column D (style index referencing the
# "YYYY-MM-DD HH:MM:SS" numFmt already applied to D2:D352).
$ws.Range("D293").NumberFormat = $ws.Range("D294").NumberFormat
